$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Edit 1: Append new sentence (with spell-check proofErr markers around the
# names "Vladmir" and "Estivill") to the paragraph that ends with
# "There are a variety of ways in which one could define how sorted a list
# is. ", then drop one of the two now-redundant blank paragraphs that
# follow it.
# -------------------------------------------------------------------------

$anchor1 = $d.Content
$found1 = $anchor1.Find.Execute(
    "There are a variety of ways in which one could define how sorted a list is. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPoint1 = $anchor1.End
$addedText1 = "The general consensus is that Vladmir Estivill-Castro and Derick Wood give the best axiomatic definitions of measures of disorder. "

$ins1 = $d.Range($insertPoint1, $insertPoint1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The general consensus is that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vladmir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Estivill</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-Castro and Derick Wood give the best axiomatic definitions of measures of disorder. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins1.InsertXML($xml1)

# InsertXML created the new text as its own paragraph; merge it back into
# the previous paragraph by removing the paragraph mark that now separates
# them.
$mark1 = $d.Range($insertPoint1, $insertPoint1 + 1)
$mark1.Delete()

$afterInsert1 = $insertPoint1 + $addedText1.Length

# Two blank paragraphs used to follow; the diff keeps only one of them.
$extraBlank1 = $d.Range($afterInsert1 + 1, $afterInsert1 + 2)
$extraBlank1.Delete()

# -------------------------------------------------------------------------
# Edit 2: Insert a new paragraph of body text between the two blank
# paragraphs that sit right after the "Methodology" heading (before the
# horizontal-rule drawing paragraph).
# -------------------------------------------------------------------------

$anchor2 = $d.Content
$found2 = $anchor2.Find.Execute(
    "Methodology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingEnd2 = $anchor2.End
$firstBlankMark2 = $headingEnd2 + 1

$ins2 = $d.Range($firstBlankMark2, $firstBlankMark2)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t xml:space="preserve">I decided to have different test class which capture the time it takes to run the tests. This was primarily for clean and modular code. It also allowed me to easily run the analysis by running the class or by using IntelliJ&#8217;s test runner. The result, which was stored in a CSV, was then put into a pandas data frame for analysis. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins2.InsertXML($xml2)

Write-Output "done"
